$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.444.20"
$ws.Range("E2").Value = "  +4.60%  "

$ws.Range("D3").Value = "'2.489.07"
$ws.Range("E3").Value = "  +2.71%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'322.79"
$ws.Range("E5").Value = "  +1.39%  "

$ws.Range("D6").Value = "'105.16"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").Value = "'0.524"
$ws.Range("E7").Value = "  +1.62%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +2.90%  "

$ws.Range("D10").Value = "'38.18"
$ws.Range("E10").Value = "  +7.37%  "

$ws.Range("D11").Value = "'0.0812"
$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("E12").Value = "  +1.06%  "

$ws.Range("D13").Value = "'18.35"
$ws.Range("E13").Value = "  +0.96%  "

$ws.Range("D14").Value = "'7.18"
$ws.Range("E14").Value = "  +1.68%  "

$ws.Range("D15").Value = "'2.879.12"
$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "'2.497.57"
$ws.Range("E16").Value = "  +2.98%  "

$ws.Range("D17").Value = "'0.847"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("D18").Value = "'47.338.12"
$ws.Range("E18").Value = "  +4.59%  "

$ws.Range("D19").Value = "'12.76"
$ws.Range("E19").Value = "  +4.57%  "

$ws.Range("D20").Value = "'6.57"
$ws.Range("E20").Value = "  +3.55%  "

$ws.Range("D21").Value = "'0.0₃0937"
$ws.Range("E21").Value = "  +1.46%  "

$ws.Range("E22").Value = "  +2.61%  "

$ws.Range("E23").Value = "  +6.32%  "

$ws.Range("D24").Value = "'251.27"
$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("D25").Value = "'2.58"
$ws.Range("E25").Value = "  +3.60%  "

$ws.Range("D26").Value = "'26.20"

$ws.Range("E27").Value = "  -0.07%  "

$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  +4.56%  "

$ws.Range("E29").Value = "  +6.43%  "

$ws.Range("D30").Value = "'35.13"
$ws.Range("E30").Value = "  +6.67%  "

$ws.Range("D31").Value = "'0.134"
$ws.Range("E31").Value = "  +8.00%  "

$ws.Range("E32").Value = "  +0.32%  "

$ws.Range("D33").Value = "'19.66"
$ws.Range("E33").Value = "  -3.17%  "

$ws.Range("D34").Value = "'5.38"
$ws.Range("E34").Value = "  +3.35%  "

$ws.Range("E35").Value = "  +2.46%  "

$ws.Range("D36").Value = "'1.01"
$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("E37").Value = "  +5.66%  "

$ws.Range("E38").Value = "  +3.96%  "

$ws.Range("E39").Value = "  +4.10%  "

$ws.Range("E40").Value = "  +2.12%  "

$ws.Range("D41").Value = "'2.24"
$ws.Range("E41").Value = "  +1.51%  "

$ws.Range("D42").Value = "'121.73"
$ws.Range("E42").Value = "  -3.25%  "

$ws.Range("D43").Value = "'21.37"
$ws.Range("E43").Value = "  +3.56%  "

$ws.Range("D44").Value = "'0.0297"
$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("D45").Value = "'1.964.55"
$ws.Range("E45").Value = "  +1.79%  "

$ws.Range("D46").Value = "'2.98"
$ws.Range("E46").Value = "  +2.17%  "

$ws.Range("E47").Value = "  -0.43%  "

$ws.Range("D48").Value = "'1.80"
$ws.Range("E48").Value = "  +1.24%  "

$ws.Range("D49").Value = "'9.16"
$ws.Range("E49").Value = "  +0.39%  "

$ws.Range("D50").Value = "'5.26"
$ws.Range("E50").Value = "  +11.56%  "

$ws.Range("D51").Value = "'79.49"
$ws.Range("E51").Value = "  +3.63%  "
